# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-row observations (Fecha, Calidad,
# Volumen, Precio minimo/maximo/promedio, Unidad de comercializacion,
# Precio $/Kg, Kg/unidad) across rows 2-19 while keeping the row's
# identity columns (Mercado, Region, Codreg, Tipo, Producto, Categoria,
# Variedad, Origen) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a unit, by column letter.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# Snapshot of the "before" values for every affected cell, keyed by
# "<col><row>".
$orig = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 19; $r++) {
        $orig["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# Target row <- source row (old row whose values now populate target row).
$map = @{
    2  = 4
    3  = 8
    4  = 10
    5  = 16
    6  = 18
    7  = 17
    8  = 6
    9  = 3
    10 = 15
    11 = 13
    12 = 19
    13 = 5
    14 = 12
    15 = 9
    16 = 2
    17 = 7
    18 = 14
    19 = 11
}

foreach ($targetRow in $map.Keys) {
    $sourceRow = $map[$targetRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $orig["$col$sourceRow"]
    }
}
